$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New due dates (Excel serial date numbers) - payments now start 3 months
# earlier and fall on the 2nd of the month instead of the 6th.
$newDates = @(44959, 44987, 45018, 45048, 45079, 45109, 45140, 45171, 45201, 45232, 45262, 45293)

# New repayment amount: the monthly payment is split/halved.
$newAmount = 45833.33333333334

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
    $ws.Cells.Item($row, 3).Value = $newAmount
}
